$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "RMS"
$ws.Range("J1").Value = "Entropy"
$ws.Range("K1").Value = "Mean Absolute Value"
$ws.Range("L1").Value = "Wilson Amplitude"
$ws.Range("M1").Value = "Zero Crossing"

$ws.Columns.Item(11).ColumnWidth = 17.608072916666668
$ws.Columns.Item(12).ColumnWidth = 14.498697916666666
$ws.Columns.Item(13).ColumnWidth = 11.053385416666666

$ws.Range("N1").Select() | Out-Null
